# templateImportOrders.xlsx -> "new version" template
#
# 1. Rename the original (only) sheet "Sheet1" -> "your data".
# 2. Add a new sheet "legend" right after it, holding the three lookup
#    lists used by the data-validation dropdowns on "your data":
#       A: PickupType   (LATER / NOW / REGULER)
#       B: Boolean      (YES / NO)
#       C: Payment Type (WALLET / CASH)
# 3. Update the two sample-row cells whose literal content changed
#    (S2: "Later" -> "REGULER", Y2: fixed a stray space before the colon).
# 4. Convert the TRUE/FALSE sample cells (U2, V2, W2, AH2) from booleans
#    (one of them formula-derived) into plain YES/NO text, matching the
#    new list-validated columns.
# 5. Wire up the list data validations against the "legend" sheet.
# 6. Restore selection/active sheet so "your data" is the visible tab.

$wb = $excel.ActiveWorkbook

# --- 1. rename existing sheet, 2. add "legend" sheet right after it ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "your data"

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "legend"

# --- legend sheet contents ---
$ws2.Range("A1").Value = "PickupType"
$ws2.Range("B1").Value = "Boolean"
$ws2.Range("A2").Value = "LATER"
$ws2.Range("A3").Value = "NOW"
$ws2.Range("A4").Value = "REGULER"
$ws2.Range("B2").Value = "YES"
$ws2.Range("B3").Value = "NO"
$ws2.Range("C1").Value = "Payment Type"
$ws2.Range("C2").Value = "WALLET"
$ws2.Range("C3").Value = "CASH"
$null = $ws2.Range("C4").Select()

# --- 3/4. sample row updates on "your data" ---
$ws1.Range("S2").Value = "REGULER"
$ws1.Range("U2").Value = "NO"
$ws1.Range("V2").Value = "NO"
$ws1.Range("W2").Value = "YES"
$ws1.Range("Y2").Value = "Isinya berupa kaos. Nomor Penerima 2: 0816100600"
$ws1.Range("AH2").Value = "NO"

# --- 5. data validation dropdowns sourced from "legend" ---
$ws1.Range("S1:S1048576").Validation.Add(3, 1, 1, "=legend!`$A`$2:`$A`$4")
$ws1.Range("U1:W1048576").Validation.Add(3, 1, 1, "=legend!`$B`$2:`$B`$3")
$ws1.Range("AH1:AH1048576").Validation.Add(3, 1, 1, "=legend!`$B`$2:`$B`$3")
$ws1.Range("AA1:AA1048576").Validation.Add(3, 1, 1, "=legend!`$C`$2:`$C`$3")

# --- 6. leave "your data" as the active/visible sheet & selection ---
$null = $ws1.Activate()
$null = $ws1.Range("AA2").Select()
